$d = $word.ActiveDocument

# 1) Paragraph 1: merge "PDF417 code that encodes the text "richardfxr."" into a single run,
#    removing the spellcheck proofErr markers around "richardfxr".
$d.Content.Find.Execute(
    "PDF417 code that encodes the text " + [char]8220 + "richardfxr." + [char]8221,
    $false, $false, $false, $false, $false, $true, 1, $false,
    "PDF417 code that encodes the text " + [char]8220 + "richardfxr." + [char]8221,
    2)

# 2) Insert the new alt-text paragraph for the About card image, using the
#    (currently empty) paragraph right after the PDF417 paragraph.
$pAlt = $d.Paragraphs.Item(2)
$pAlt.Range.Text = "Graphite self-portrait. I" + [char]8217 + "m an Asian male with short black hair wearing a pair of aviator-style glasses. "

# 3) Last paragraph: merge "Not available for projects in the near future." into a
#    single run, removing the grammar-check proofErr markers.
$d.Content.Find.Execute(
    "Not available for projects in the near future.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Not available for projects in the near future.",
    2)

Write-Output "done"
